$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New corrected values for the "Total" (B) and "Community" (D) columns,
# rows 2-13 (months 1-12), per the commit "Fixing the big mistake".
$values = @(
    @{ Row = 2;  B = 7873.70901350002;   D = 505.2825002833334 },
    @{ Row = 3;  B = 7366.369618183351;  D = 479.8452055333333 },
    @{ Row = 4;  B = 7900.730378083354;  D = 529.6545940666667 },
    @{ Row = 5;  B = 7617.880432633352;  D = 495.8968661333333 },
    @{ Row = 6;  B = 7903.597525466686;  D = 523.54358935 },
    @{ Row = 7;  B = 7655.652308883353;  D = 509.1682569166667 },
    @{ Row = 8;  B = 7895.98543095002;   D = 518.8358951499999 },
    @{ Row = 9;  B = 7895.23621675002;   D = 523.1529089666667 },
    @{ Row = 10; B = 7664.542564450019;  D = 500.90419505 },
    @{ Row = 11; B = 7893.204864216686;  D = 523.8685568833333 },
    @{ Row = 12; B = 7666.541936400019;  D = 492.1805045666667 },
    @{ Row = 13; B = 7641.802333766685;  D = 504.9592266333333 }
)

foreach ($entry in $values) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
    $ws.Cells.Item($entry.Row, 4).Value = $entry.D
}
